$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the next data row (row 19) of the "Avances Etiquetado Roboflow" log
# with the 5/20/2025 update: same counts as the previous entry, new date and
# a note from Rafael about finals getting in the way of labeling.
$ws.Range("D19").Value = "20/5/2025"
$ws.Range("E19").Value = 135
$ws.Range("F19").Value = 218
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 650
$ws.Range("J19").Value = "no he podido etiquetar por examenes finales una disculpa (rafael)"

# Give the table a fresh blank row (row 20) to keep logging future entries,
# matching the formatting of the row that just got filled in.
$ws.Range("D19:J19").Copy($ws.Range("D20:J20"))
$ws.Range("D20:J20").ClearContents()
$ws.Rows.Item(20).RowHeight = $ws.Rows.Item(19).RowHeight

# Grow Table1 so the new blank row is included in the table range.
$t = $ws.ListObjects.Item(1)
$t.Resize($ws.Range("D4:J20"))

# Reflect where the user's selection ended up after adding the new row.
$ws.Range("J23").Select()
